# Project DesignFirst save: update rule R30's "Integer min" threshold.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell C10 ("Integer min" for rule R30) changes from 18 to 100.
$ws.Range("C10").Value = 100
